$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values (new randomized test data for TC_001)
$ws.Range("A1").Value = "Bangaloreeez"
$ws.Range("A2").Value = "Bangaloreyyq"
$ws.Range("B1").Value = "Lifeuuv"
$ws.Range("B2").Value = "Lifeiima"
$ws.Range("C2").Value = "'9764454332"
$ws.Range("C1").Value = "'9341262932"

# Move active selection from D6 to C4
$ws.Range("C4").Select()
